$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3 ---
$ws.Range("G3").Value = 0
$ws.Range("N3").Value = 0

# --- Row 4 ---
$ws.Range("G4").Value = 1
$ws.Range("N4").Value = 1

# --- Row 6 ---
$ws.Range("G6").Value = 0
$ws.Range("N6").Value = 0

# --- Row 7 ---
$ws.Range("C7").Value = 1143
$ws.Range("G7").Value = 0
$ws.Range("J7").Value = 1354
$ws.Range("N7").Value = 0

# F7 / M7 hold text-like numeric strings (shared strings), not real numbers.
# Prefix with an apostrophe so Excel stores them as text, then restore the
# default "Normal" style so no stray per-cell number format lingers.
$ws.Range("F7").Value = "'0.48846153846153845"
$ws.Range("M7").Value = "'0.5786324786324787"
$ws.Range("F7").Style = "Normal"
$ws.Range("M7").Style = "Normal"
